$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 685.5714
$ws.Range("I2").Value = 685.5714
$ws.Range("K2").Value = 685.5714
$ws.Range("M2").Value = -572.5714
$ws.Range("H28").Value = 1314.5
$ws.Range("J28").Value = 2475.8
$ws.Range("L28").Value = 2475.8
$ws.Range("N28").Value = -3445.8
$ws.Range("H51").Value = 12499.5
$ws.Range("I51").Value = 9999
$ws.Range("J51").Value = 13333
$ws.Range("K51").Value = 9999
$ws.Range("L51").Value = 13333
$ws.Range("M51").Value = -9515
$ws.Range("N51").Value = -14301
$ws.Range("H70").Value = 2659.4
$ws.Range("I70").Value = 2324.25
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 6972.75
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -6702.75
$ws.Range("N70").Value = -12540
$ws.Range("H73").Value = 2659.4
$ws.Range("I73").Value = 2324.25
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 6972.75
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = -6036.75
$ws.Range("N73").Value = -13872
$ws.Range("H92").Value = 52631940
$ws.Range("J92").Value = 530.4
$ws.Range("L92").Value = 530.4
$ws.Range("N92").Value = -3026.4
$ws.Range("H93").Value = 43344.5
$ws.Range("J93").Value = 43344.5
$ws.Range("L93").Value = 43344.5
$ws.Range("N93").Value = -48336.5
$ws.Range("H99").Value = 111111200
$ws.Range("I99").Value = 111111200
$ws.Range("K99").Value = 333333600
$ws.Range("M99").Value = -333332102
$ws.Range("H111").Value = 2309.25
$ws.Range("I111").Value = 3251.7144
$ws.Range("J111").Value = 989.8
$ws.Range("K111").Value = 9755.143199999999
$ws.Range("L111").Value = 2969.4
$ws.Range("M111").Value = -6688.143199999999
$ws.Range("N111").Value = -9103.4
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("H127").Value = 865.3333
$ws.Range("I127").Value = 938.4
$ws.Range("K127").Value = 2815.2
$ws.Range("M127").Value = 2144.8
$ws.Range("H132").Value = 42879.312
$ws.Range("I132").Value = 47933.07
$ws.Range("K132").Value = 143799.21
$ws.Range("M132").Value = -141269.21
$ws.Range("H137").Value = 1582.8334
$ws.Range("I137").Value = 974.75
$ws.Range("J137").Value = 2799
$ws.Range("K137").Value = 2924.25
$ws.Range("L137").Value = 8397
$ws.Range("M137").Value = -374.25
$ws.Range("N137").Value = -13497
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 4400
$ws.Range("I15").Value = 3800
$ws.Range("J15").Value = 5000
$ws.Range("K15").Value = 3800
$ws.Range("L15").Value = 5000
$ws.Range("M15").Value = -3450
$ws.Range("N15").Value = -5700
$ws.Range("H32").Value = 1470.4166
$ws.Range("I32").Value = 967.7273
$ws.Range("J32").Value = 7000
$ws.Range("K32").Value = 967.7273
$ws.Range("L32").Value = 7000
$ws.Range("M32").Value = -680.7273
$ws.Range("N32").Value = -7574
$ws.Range("H106").Value = 9999
$ws.Range("J106").Value = 9999
$ws.Range("L106").Value = 9999
$ws.Range("N106").Value = -12523
$ws.Range("H132").Value = 3253.25
$ws.Range("I132").Value = 2999.6667
$ws.Range("J132").Value = 4014
$ws.Range("K132").Value = 8999.000100000001
$ws.Range("L132").Value = 12042
$ws.Range("M132").Value = -6469.000100000001
$ws.Range("N132").Value = -17102

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1436.75
$ws.Range("I80").Value = 1436.5
$ws.Range("K80").Value = 1436.5
$ws.Range("M80").Value = -438.5
$ws.Range("H83").Value = 1436.75
$ws.Range("I83").Value = 1436.5
$ws.Range("K83").Value = 7182.5
$ws.Range("M83").Value = -2190.5
$ws.Range("H102").Value = 16732.3
$ws.Range("J102").Value = 49612
$ws.Range("L102").Value = 49612
$ws.Range("N102").Value = -56102

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 23599.6
$ws.Range("J92").Value = 23599.6
$ws.Range("L92").Value = 23599.6
$ws.Range("N92").Value = -28591.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 51.333332
$ws.Range("I33").Value = 47
$ws.Range("J33").Value = 60
$ws.Range("K33").Value = 282
$ws.Range("L33").Value = 360
$ws.Range("M33").Value = 1
$ws.Range("N33").Value = -926
$ws.Range("H68").Value = 2993.2727
$ws.Range("I68").Value = 2993.5
$ws.Range("K68").Value = 8980.5
$ws.Range("M68").Value = -8169.5
$ws.Range("H71").Value = 2993.2727
$ws.Range("I71").Value = 2993.5
$ws.Range("K71").Value = 26941.5
$ws.Range("M71").Value = -22885.5
$ws.Range("H92").Value = 692.5
$ws.Range("I92").Value = 690
$ws.Range("J92").Value = 700
$ws.Range("K92").Value = 2070
$ws.Range("L92").Value = 2100
$ws.Range("M92").Value = -822
$ws.Range("N92").Value = -4596
$ws.Range("H113").Value = 687.4
$ws.Range("J113").Value = 467
$ws.Range("L113").Value = 1401
$ws.Range("N113").Value = -5741
$ws.Range("H121").Value = 464.66666
$ws.Range("I121").Value = 282.66666
$ws.Range("K121").Value = 847.9999799999999
$ws.Range("M121").Value = 462.0000200000001
$ws.Range("H129").Value = 1110.25
$ws.Range("I129").Value = 633.25
$ws.Range("J129").Value = 1587.25
$ws.Range("K129").Value = 1899.75
$ws.Range("L129").Value = 4761.75
$ws.Range("M129").Value = 3100.25
$ws.Range("N129").Value = -14761.75
$ws.Range("H139").Value = 3254
$ws.Range("I139").Value = 3254
$ws.Range("K139").Value = 9762
$ws.Range("M139").Value = -4622

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 3023.4443
$ws.Range("I22").Value = 642.8
$ws.Range("J22").Value = 5999.25
$ws.Range("K22").Value = 642.8
$ws.Range("L22").Value = 5999.25
$ws.Range("M22").Value = -113.8
$ws.Range("N22").Value = -7057.25
$ws.Range("H45").Value = 36666.668
$ws.Range("J45").Value = 35000
$ws.Range("L45").Value = 35000
$ws.Range("N45").Value = -36118
$ws.Range("H70").Value = 31254682
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 33338062
$ws.Range("K70").Value = 4000
$ws.Range("L70").Value = 33338062
$ws.Range("M70").Value = -3730
$ws.Range("N70").Value = -33338602
$ws.Range("H73").Value = 31254682
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 33338062
$ws.Range("K73").Value = 4000
$ws.Range("L73").Value = 33338062
$ws.Range("M73").Value = -3064
$ws.Range("N73").Value = -33339934
$ws.Range("H122").Value = 5635.143
$ws.Range("I122").Value = 5490.3335
$ws.Range("K122").Value = 16471.0005
$ws.Range("M122").Value = -14021.0005
$ws.Range("H126").Value = 1999.25
$ws.Range("I126").Value = 1999.25
$ws.Range("K126").Value = 5997.75
$ws.Range("M126").Value = -3527.75
$ws.Range("H132").Value = 2312.6667
$ws.Range("I132").Value = 1462
$ws.Range("K132").Value = 4386
$ws.Range("M132").Value = -1856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("H103").Value = 20268.334
$ws.Range("J103").Value = 20268.334
$ws.Range("L103").Value = 20268.334
$ws.Range("N103").Value = -22612.334
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("H132").Value = 2246.75
$ws.Range("I132").Value = 1151.3334
$ws.Range("K132").Value = 3454.0002
$ws.Range("M132").Value = -924.0001999999999
$ws.Range("M16").ClearContents()
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 175666.67
$ws.Range("J15").Value = 13500
$ws.Range("L15").Value = 13500
$ws.Range("N15").Value = -14076
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("H132").Value = 1064.5
$ws.Range("I132").Value = 1064.5
$ws.Range("K132").Value = 1064.5
$ws.Range("M132").Value = -663.5
$ws.Range("N22").ClearContents()
